$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5068
$ws.Range("I32").Value = 4493
$ws.Range("K32").Value = 4493
$ws.Range("M32").Value = -4167

$ws.Range("H43").Value = 2528.6
$ws.Range("I43").Value = 3599.6
$ws.Range("J43").Value = 1457.6
$ws.Range("K43").Value = 3599.6
$ws.Range("L43").Value = 1457.6
$ws.Range("M43").Value = -3530.6
$ws.Range("N43").Value = -1595.6

$ws.Range("H51").Value = 3500
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H55").Value = 141.26315
$ws.Range("I55").Value = 123.125
$ws.Range("J55").Value = 154.45454
$ws.Range("K55").Value = 123.125
$ws.Range("L55").Value = 154.45454
$ws.Range("M55").Value = 90.875
$ws.Range("N55").Value = -582.45454

$ws.Range("H62").Value = 4581
$ws.Range("I62").Value = 3260.3076
$ws.Range("K62").Value = 3260.3076
$ws.Range("M62").Value = -2636.3076

$ws.Range("H65").Value = 4581
$ws.Range("I65").Value = 3260.3076
$ws.Range("K65").Value = 16301.538
$ws.Range("M65").Value = -13181.538

$ws.Range("H116").Value = 2888.1667
$ws.Range("J116").Value = 3516.6667
$ws.Range("L116").Value = 3516.6667
$ws.Range("N116").Value = -10400.6667

$ws.Range("H131").Value = 7499.875
$ws.Range("I131").Value = 5999.5
$ws.Range("K131").Value = 17998.5
$ws.Range("M131").Value = -12958.5

$ws.Range("H138").Value = 7579960
$ws.Range("J138").Value = 13895712
$ws.Range("L138").Value = 41687136
$ws.Range("N138").Value = -41697416

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1649.25
$ws.Range("I45").Value = 1799
$ws.Range("K45").Value = 1799
$ws.Range("M45").Value = -1422

$ws.Range("H51").Value = 38495
$ws.Range("J51").Value = 38495
$ws.Range("L51").Value = 38495
$ws.Range("N51").Value = -40007

$ws.Range("H102").Value = 3831.5
$ws.Range("I102").Value = 2183.5557
$ws.Range("J102").Value = 6797.8
$ws.Range("K102").Value = 2183.5557
$ws.Range("L102").Value = 6797.8
$ws.Range("M102").Value = -561.5556999999999
$ws.Range("N102").Value = -10041.8

$ws.Range("H110").Value = 12688.781
$ws.Range("I110").Value = 13729.5
$ws.Range("J110").Value = 5403.75
$ws.Range("K110").Value = 13729.5
$ws.Range("L110").Value = 5403.75
$ws.Range("M110").Value = -11684.5
$ws.Range("N110").Value = -9493.75

$ws.Range("H122").Value = 3082.4333
$ws.Range("I122").Value = 2110.2104
$ws.Range("K122").Value = 6330.6312
$ws.Range("M122").Value = -3880.6312

$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13745.692
$ws.Range("I86").Value = 9290.538
$ws.Range("J86").Value = 18200.846
$ws.Range("K86").Value = 9290.538
$ws.Range("L86").Value = 18200.846
$ws.Range("M86").Value = -8167.538
$ws.Range("N86").Value = -20446.846

$ws.Range("H89").Value = 13745.692
$ws.Range("I89").Value = 9290.538
$ws.Range("J89").Value = 18200.846
$ws.Range("K89").Value = 46452.69
$ws.Range("L89").Value = 91004.23000000001
$ws.Range("M89").Value = -40836.69
$ws.Range("N89").Value = -102236.23

$ws.Range("H94").Value = 1134.3
$ws.Range("I94").Value = 903.6087
$ws.Range("K94").Value = 903.6087
$ws.Range("M94").Value = -452.6087

$ws.Range("H105").Value = 6862.15
$ws.Range("I105").Value = 14236.125
$ws.Range("J105").Value = 1946.1666
$ws.Range("K105").Value = 14236.125
$ws.Range("L105").Value = 1946.1666
$ws.Range("M105").Value = -12489.125
$ws.Range("N105").Value = -5440.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4430.346
$ws.Range("I22").Value = 5343.2
$ws.Range("K22").Value = 5343.2
$ws.Range("M22").Value = -4993.2

$ws.Range("H99").Value = 9469.700000000001
$ws.Range("I99").Value = 8742.357
$ws.Range("K99").Value = 8742.357
$ws.Range("M99").Value = -7244.357

$ws.Range("H105").Value = 10993.357
$ws.Range("I105").Value = 2500.9
$ws.Range("K105").Value = 2500.9
$ws.Range("M105").Value = -753.9000000000001

$ws.Range("H107").Value = 2773.75
$ws.Range("I107").Value = 2033.3334
$ws.Range("K107").Value = 2033.3334
$ws.Range("M107").Value = -113.3334

$ws.Range("H124").Value = 75992.28999999999
$ws.Range("J124").Value = 75992.28999999999
$ws.Range("L124").Value = 75992.28999999999
$ws.Range("N124").Value = -80902.28999999999

$ws.Range("H126").Value = 9469.700000000001
$ws.Range("I126").Value = 8742.357
$ws.Range("K126").Value = 26227.071
$ws.Range("M126").Value = -23757.071

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1124.6428
$ws.Range("J5").Value = 2500
$ws.Range("L5").Value = 7500
$ws.Range("N5").Value = -7724

$ws.Range("H37").Value = 139494.25
$ws.Range("J37").Value = 139494.25
$ws.Range("L37").Value = 418482.75
$ws.Range("N37").Value = -418706.75

$ws.Range("H97").Value = 399
$ws.Range("J97").Value = 399
$ws.Range("L97").Value = 1197
$ws.Range("N97").Value = -2189

$ws.Range("H98").Value = 2309
$ws.Range("I98").Value = 642.6667
$ws.Range("J98").Value = 3308.8
$ws.Range("K98").Value = 1928.0001
$ws.Range("L98").Value = 9926.400000000001
$ws.Range("M98").Value = -430.0001
$ws.Range("N98").Value = -12922.4

$ws.Range("H133").Value = 12636.944
$ws.Range("I133").Value = 6158
$ws.Range("K133").Value = 18474
$ws.Range("M133").Value = -13414

$ws.Range("H134").Value = 7742.304
$ws.Range("I134").Value = 2071.6
$ws.Range("J134").Value = 18374.875
$ws.Range("K134").Value = 6214.799999999999
$ws.Range("L134").Value = 55124.625
$ws.Range("M134").Value = -1144.799999999999
$ws.Range("N134").Value = -65264.625

$ws.Range("H135").Value = 1124.6428
$ws.Range("J135").Value = 2500
$ws.Range("L135").Value = 22500
$ws.Range("N135").Value = -27570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4458.154
$ws.Range("J70").Value = 4999.25
$ws.Range("L70").Value = 4999.25
$ws.Range("N70").Value = -5539.25

$ws.Range("H73").Value = 4458.154
$ws.Range("J73").Value = 4999.25
$ws.Range("L73").Value = 4999.25
$ws.Range("N73").Value = -6871.25

$ws.Range("H88").Value = 62497.5
$ws.Range("J88").Value = 59495
$ws.Range("L88").Value = 59495
$ws.Range("N88").Value = -60397

$ws.Range("H91").Value = 62497.5
$ws.Range("J91").Value = 59495
$ws.Range("L91").Value = 59495
$ws.Range("N91").Value = -62615

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 394.5
$ws.Range("I16").Value = 388.9
$ws.Range("K16").Value = 388.9
$ws.Range("M16").Value = -218.9

$ws.Range("H22").Value = 3181.2104
$ws.Range("I22").Value = 1900
$ws.Range("K22").Value = 1900
$ws.Range("M22").Value = -1605

$ws.Range("H27").Value = 3181.2104
$ws.Range("I27").Value = 1900
$ws.Range("K27").Value = 1900
$ws.Range("M27").Value = -1793

$ws.Range("H40").Value = 4280.4287
$ws.Range("I40").Value = 4280.4287
$ws.Range("K40").Value = 4280.4287
$ws.Range("M40").Value = -4144.4287

$ws.Range("H46").Value = 1819.6511
$ws.Range("I46").Value = 626.6667
$ws.Range("K46").Value = 626.6667
$ws.Range("M46").Value = -438.6667

$ws.Range("H61").Value = 7166.6665
$ws.Range("I61").Value = 7500
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 7500
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -7298
$ws.Range("N61").Value = -7404

$ws.Range("H113").Value = 7166.6665
$ws.Range("I113").Value = 7500
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 7500
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = -5330
$ws.Range("N113").Value = -11340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 64249.75
$ws.Range("I57").Value = 58500
$ws.Range("J57").Value = 69999.5
$ws.Range("K57").Value = 58500
$ws.Range("L57").Value = 69999.5
$ws.Range("M57").Value = -57746
$ws.Range("N57").Value = -71507.5

$ws.Range("H62").Value = 7385.5454
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 7674.1
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 7674.1
$ws.Range("M62").Value = -3876
$ws.Range("N62").Value = -8922.1

$ws.Range("H65").Value = 7385.5454
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 7674.1
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 38370.5
$ws.Range("M65").Value = -19380
$ws.Range("N65").Value = -44610.5

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
